$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 13.14763442785007
$ws.Range("C2").Value = 5.67433876543938
$ws.Range("D2").Value = 6.99581215405264
$ws.Range("E2").Value = 10.89328524447625
$ws.Range("F2").Value = 36.54564082404067
$ws.Range("K2").Value = 12.54395826449712
$ws.Range("M2").Value = 15.45836519353906
$ws.Range("N2").Value = 22.05602951071927
$ws.Range("B3").Value = 12.90300152081766
$ws.Range("C3").Value = 5.468126719357647
$ws.Range("D3").Value = 7.004714139653095
$ws.Range("E3").Value = 10.68004866482388
$ws.Range("F3").Value = 36.32024241319999
$ws.Range("K3").Value = 12.38199045569924
$ws.Range("M3").Value = 15.30570589124348
$ws.Range("N3").Value = 22.0791869710074
$ws.Range("B4").Value = 12.75462639527463
$ws.Range("C4").Value = 5.339488557428746
$ws.Range("D4").Value = 7.010280847802457
$ws.Range("E4").Value = 10.55058114629476
$ws.Range("F4").Value = 36.18938287075332
$ws.Range("K4").Value = 12.28525165783887
$ws.Range("M4").Value = 15.21569009320452
$ws.Range("N4").Value = 22.09520667110874
$ws.Range("B5").Value = 12.69471338035317
$ws.Range("C5").Value = 5.286654779910189
$ws.Range("D5").Value = 7.012574888963263
$ws.Range("E5").Value = 10.4982665280688
$ws.Range("F5").Value = 36.13798917063346
$ws.Range("K5").Value = 12.24656038395238
$ws.Range("M5").Value = 15.1799809145309
$ws.Range("N5").Value = 22.10218718289541
$ws.Range("B6").Value = 12.68480080359088
$ws.Range("C6").Value = 5.277859812751478
$ws.Range("D6").Value = 7.012957362538873
$ws.Range("E6").Value = 10.48960882402869
$ws.Range("F6").Value = 36.12957304756893
$ws.Range("K6").Value = 12.24018125590423
$ws.Range("M6").Value = 15.17411127187028
$ws.Range("N6").Value = 22.10337359759253
$ws.Range("B7").Value = 12.75381603632778
$ws.Range("C7").Value = 5.338777560639899
$ws.Range("D7").Value = 7.010311682299976
$ws.Range("E7").Value = 10.54987371151465
$ws.Range("F7").Value = 36.18868188506234
$ws.Range("K7").Value = 12.28472683205454
$ws.Range("M7").Value = 15.21520451832633
$ws.Range("N7").Value = 22.09529898191576
$ws.Range("B8").Value = 13.06295948568621
$ws.Range("C8").Value = 5.603711838228573
$ws.Range("D8").Value = 6.998860750247013
$ws.Range("E8").Value = 10.81950435178484
$ws.Range("F8").Value = 36.46638290454345
$ws.Range("K8").Value = 12.48757919285919
$ws.Range("M8").Value = 15.404980202063
$ws.Range("N8").Value = 22.06364012056168
$ws.Range("B9").Value = 13.6797506061947
$ws.Range("C9").Value = 6.103395136797225
$ws.Range("D9").Value = 6.977195878263916
$ws.Range("E9").Value = 11.35647424617686
$ws.Range("F9").Value = 37.06894966650871
$ws.Range("K9").Value = 12.90465294158327
$ws.Range("M9").Value = 15.80485458249608
$ws.Range("N9").Value = 22.01586984991548
$ws.Range("B10").Value = 14.13414272667914
$ws.Range("C10").Value = 6.453910265604182
$ws.Range("D10").Value = 6.961745832631292
$ws.Range("E10").Value = 11.75159115121181
$ws.Range("F10").Value = 37.54440972837771
$ws.Range("K10").Value = 13.21979691446713
$ws.Range("M10").Value = 16.11305724218503
$ws.Range("N10").Value = 21.98953096584894
$ws.Range("B11").Value = 14.34007259902635
$ws.Range("C11").Value = 6.608974829726868
$ws.Range("D11").Value = 6.954815415646053
$ws.Range("E11").Value = 11.93058177620859
$ws.Range("F11").Value = 37.76721417899348
$ws.Range("K11").Value = 13.36440868978314
$ws.Range("M11").Value = 16.25584421256922
$ws.Range("N11").Value = 21.97945724463403
$ws.Range("B12").Value = 14.41785822520785
$ws.Range("C12").Value = 6.667007435208581
$ws.Range("D12").Value = 6.952204887519605
$ws.Range("E12").Value = 11.9981829517462
$ws.Range("F12").Value = 37.85246691022427
$ws.Range("K12").Value = 13.4192956737165
$ws.Range("M12").Value = 16.31023829092706
$ws.Range("N12").Value = 21.97591746928709
$ws.Range("B13").Value = 14.40111569064711
$ws.Range("C13").Value = 6.654540450158104
$ws.Range("D13").Value = 6.952766497845801
$ws.Range("E13").Value = 11.983632865834
$ws.Range("F13").Value = 37.83406790290271
$ws.Range("K13").Value = 13.40747004538657
$ws.Range("M13").Value = 16.29850991687966
$ws.Range("N13").Value = 21.97666758715931
$ws.Range("B14").Value = 14.34647645529744
$ws.Range("C14").Value = 6.61376327244238
$ws.Range("D14").Value = 6.954600369124492
$ws.Range("E14").Value = 11.93614733975233
$ws.Range("F14").Value = 37.77421062872057
$ws.Range("K14").Value = 13.36892206433737
$ws.Range("M14").Value = 16.26031303211667
$ws.Range("N14").Value = 21.97916051215162
$ws.Range("B15").Value = 14.31298045481057
$ws.Range("C15").Value = 6.58869505951459
$ws.Range("D15").Value = 6.955725468449484
$ws.Range("E15").Value = 11.90703577216951
$ws.Range("F15").Value = 37.73765944736945
$ws.Range("K15").Value = 13.34532508085126
$ws.Range("M15").Value = 16.23695707720405
$ws.Range("N15").Value = 21.98072332026368
$ws.Range("B16").Value = 14.12066184636207
$ws.Range("C16").Value = 6.443683198650641
$ws.Range("D16").Value = 6.962200703847016
$ws.Range("E16").Value = 11.73987246422137
$ws.Range("F16").Value = 37.52997534500951
$ws.Range("K16").Value = 13.21036671152926
$ws.Range("M16").Value = 16.1037736911657
$ws.Range("N16").Value = 21.99022774876607
$ws.Range("B17").Value = 14.00242095687737
$ws.Range("C17").Value = 6.353557385950028
$ws.Range("D17").Value = 6.966197976021695
$ws.Range("E17").Value = 11.63707957538851
$ws.Range("F17").Value = 37.40419777405374
$ws.Range("K17").Value = 13.12785569017456
$ws.Range("M17").Value = 16.02269787524382
$ws.Range("N17").Value = 21.99654748895103
$ws.Range("B18").Value = 13.93434372916203
$ws.Range("C18").Value = 6.301309415409388
$ws.Range("D18").Value = 6.968506327115033
$ws.Range("E18").Value = 11.57788925174076
$ws.Range("F18").Value = 37.33247052824365
$ws.Range("K18").Value = 13.08051850428395
$ws.Range("M18").Value = 15.97631131828413
$ws.Range("N18").Value = 22.0003619911175
$ws.Range("B19").Value = 13.91128486317234
$ws.Range("C19").Value = 6.283550562876834
$ws.Range("D19").Value = 6.969289485923774
$ws.Range("E19").Value = 11.55783920730951
$ws.Range("F19").Value = 37.3082924893812
$ws.Range("K19").Value = 13.06451339186437
$ws.Range("M19").Value = 15.9606493501595
$ws.Range("N19").Value = 22.00168433729714
$ws.Range("B20").Value = 14.01501558105701
$ws.Range("C20").Value = 6.363194293991879
$ws.Range("D20").Value = 6.96577150625796
$ws.Range("E20").Value = 11.64802946828524
$ws.Range("F20").Value = 37.41752358228686
$ws.Range("K20").Value = 13.13662700537894
$ws.Range("M20").Value = 16.03130340894724
$ws.Range("N20").Value = 21.9958561545944
$ws.Range("B21").Value = 14.36253127990602
$ws.Range("C21").Value = 6.625759575905247
$ws.Range("D21").Value = 6.954061341878909
$ws.Range("E21").Value = 11.95010038609675
$ws.Range("F21").Value = 37.79176868658884
$ws.Range("K21").Value = 13.38024155859621
$ws.Range("M21").Value = 16.27152397840816
$ws.Range("N21").Value = 21.97842081390433
$ws.Range("B22").Value = 14.58847703930443
$ws.Range("C22").Value = 6.793335638958384
$ws.Range("D22").Value = 6.946488823734175
$ws.Range("E22").Value = 12.14644939430518
$ws.Range("F22").Value = 38.04147085826349
$ws.Range("K22").Value = 13.54016567869133
$ws.Range("M22").Value = 16.43038686494758
$ws.Range("N22").Value = 21.96862852710805
$ws.Range("B23").Value = 14.46801982303411
$ws.Range("C23").Value = 6.704282274616155
$ws.Range("D23").Value = 6.95052309855204
$ws.Range("E23").Value = 12.04177476212757
$ws.Range("F23").Value = 37.90775113728061
$ws.Range("K23").Value = 13.45476405442337
$ws.Range("M23").Value = 16.34544393936205
$ws.Range("N23").Value = 21.97370802562486
$ws.Range("B24").Value = 14.00932185327426
$ws.Range("C24").Value = 6.358838793610612
$ws.Range("D24").Value = 6.965964281225674
$ws.Range("E24").Value = 11.64307930699885
$ws.Range("F24").Value = 37.41149716117737
$ws.Range("K24").Value = 13.13266117970096
$ws.Range("M24").Value = 16.02741214265777
$ws.Range("N24").Value = 21.99616814243786
$ws.Range("B25").Value = 13.51230730683013
$ws.Range("C25").Value = 5.970839404178731
$ws.Range("D25").Value = 6.982973701827111
$ws.Range("E25").Value = 11.21079204753005
$ws.Range("F25").Value = 36.89998846239448
$ws.Range("K25").Value = 12.79007103736914
$ws.Range("M25").Value = 15.6939614876504
$ws.Range("N25").Value = 22.02725728114231
